# [ADD] Add robustness experiments
#
# Renames the existing "Paolo" user to "Pöl" (short form) on the PD/PID
# controller step/dynamic-response rows, and fills in the two previously
# blank rows (experiment 15 and 16 -> sheet rows 18 and 19) with new
# "robustness" experiments, authored by user "Cere".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New rows 18 (PD) & 19 (PID): robustness experiments, filled column by column ---
$ws.Range("C18").Value = "PD controller test robustness"
$ws.Range("C19").Value = "PID controller test robustness"

$ws.Range("D18").Value = "Control the system via the PD controller and set the theta`nreference to a constant value, when stable, slightly hit the pendulum to generate some oscillations"
$ws.Range("D19").Value = "Control the system via the PID controller and set the theta`nreference to a constant value, when stable, slightly hit the pendulum to generate some oscillations"

$ws.Range("E18").Value = "Check robustness of the controller"
$ws.Range("E19").Value = "Check robustness of the controller"

$ws.Range("G18").Value = "Cere"
$ws.Range("G19").Value = "Cere"

# --- Rename existing user "Paolo" -> "Pöl" on rows 14-17 ---
$ws.Range("G14").Value = "Pöl"
$ws.Range("G15").Value = "Pöl"
$ws.Range("G16").Value = "Pöl"
$ws.Range("G17").Value = "Pöl"

# --- Leave the selection where the editor last left off ---
$ws.Range("F12").Select() | Out-Null
